$wb = $excel.ActiveWorkbook

# ---- sheet7 (保險) ----
$ws = $wb.Worksheets.Item("保險")
$ws.Range("D1").Copy()
$ws.Range("E1:K1").PasteSpecial(-4122)
$ws.Range("D2:D4").Copy()
$ws.Range("E2:K4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"
$ws.Range("A2").Value = 91
$ws.Range("A3").Value = 92
$ws.Range("A4").Value = 93
$ws.Range("B2").Value = "富邦人壽"
$ws.Range("B3").Value = "富邦人壽"
$ws.Range("B4").Value = "中國人壽"
$ws.Range("C2").Value = "吉祥變額萬能终身壽險(A型）"
$ws.Range("C3").Value = "吉祥變額萬能終身壽險(A型）"
$ws.Range("C4").Value = "限期繳k單利增值終身壽險已型"
$ws.Range("D2").Value = "潘维剛"
$ws.Range("D3").Value = "潘維剛"
$ws.Range("D4").Value = "田正超"
$ws.Range("E2").Value = "insurance"
$ws.Range("E3").Value = "insurance"
$ws.Range("E4").Value = "insurance"
$ws.Range("F2").Value = "normal"
$ws.Range("F3").Value = "normal"
$ws.Range("F4").Value = "normal"
$ws.Range("G2").Value = "2011-11-21"
$ws.Range("G3").Value = "2011-11-21"
$ws.Range("G4").Value = "2011-11-21"
$ws.Range("H2").Value = "潘維剛"
$ws.Range("H3").Value = "潘維剛"
$ws.Range("H4").Value = "潘維剛"
$ws.Range("I2").Value = 678
$ws.Range("I3").Value = 678
$ws.Range("I4").Value = 678
$ws.Range("J2").Value = "tmpcafb1"
$ws.Range("J3").Value = "tmpcafb1"
$ws.Range("J4").Value = "tmpcafb1"
$ws.Range("K2").Value = 91
$ws.Range("K3").Value = 92
$ws.Range("K4").Value = 93

# ---- sheet8 (債務) ----
$ws = $wb.Worksheets.Item("債務")
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$ws.Range("G2:G3").Copy()
$ws.Range("H2:N3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"
$ws.Range("A2").Value = 103
$ws.Range("A3").Value = 104
$ws.Range("B2").Value = "—皞借款"
$ws.Range("B3").Value = "一般借款"
$ws.Range("C2").Value = "田正超"
$ws.Range("C3").Value = "田正超"
$ws.Range("D2").Value = "永骂銀行西松分拧臺北市松山區東興路"
$ws.Range("D3").Value = "永豐銀行板新分行新北市板橋區民權路"
$ws.Range("E2").Value = 11798308
$ws.Range("E3").Value = 7355323
$ws.Range("F2").Value = "95年10月27B"
$ws.Range("F3").Value = "96年07月25日"
$ws.Range("G2").Value = "房貸"
$ws.Range("G3").Value = "房貸"
$ws.Range("H2").Value = "debt"
$ws.Range("H3").Value = "debt"
$ws.Range("I2").Value = "normal"
$ws.Range("I3").Value = "normal"
$ws.Range("J2").Value = "2011-11-21"
$ws.Range("J3").Value = "2011-11-21"
$ws.Range("K2").Value = "潘維剛"
$ws.Range("K3").Value = "潘維剛"
$ws.Range("L2").Value = 678
$ws.Range("L3").Value = 678
$ws.Range("M2").Value = "tmpcafb1"
$ws.Range("M3").Value = "tmpcafb1"
$ws.Range("N2").Value = 103
$ws.Range("N3").Value = 104

# ---- sheet9 (事業投資) ----
$ws = $wb.Worksheets.Item("事業投資")
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$ws.Range("G2:G5").Copy()
$ws.Range("H2:N5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "owner"
$ws.Range("C1").Value = "company"
$ws.Range("D1").Value = "address"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"
$ws.Range("A2").Value = 109
$ws.Range("A3").Value = 110
$ws.Range("A4").Value = 111
$ws.Range("A5").Value = 112
$ws.Range("B2").Value = "潘維剛"
$ws.Range("B3").Value = "潘維剛"
$ws.Range("B4").Value = "田正超"
$ws.Range("B5").Value = "田正超"
$ws.Range("C2").Value = "傳智國際文化事業座份有限公司"
$ws.Range("C3").Value = "女人網股份有限公司"
$ws.Range("C4").Value = "春田國際有限公司"
$ws.Range("C5").Value = "中餘國際股份有限公司"
$ws.Range("D2").Value = "臺北市羅斯福路2段116號3樓"
$ws.Range("D3").Value = "臺北市杭州南路1段63號5樓之1"
$ws.Range("D4").Value = "臺北市民權東路3段106巷36號7樓"
$ws.Range("D5").Value = "新北市寶橋路235巷65號5樓"
$ws.Range("E2").Value = 2000000
$ws.Range("E3").Value = 900000
$ws.Range("E4").Value = 5000000
$ws.Range("E5").Value = 612900
$ws.Range("F2").Value = "89年08月11曰"
$ws.Range("F3").Value = "98年08月10日"
$ws.Range("F4").Value = "92年03月20日"
$ws.Range("F5").Value = "82年06月28日"
$ws.Range("G2").Value = "投資"
$ws.Range("G3").Value = "投資"
$ws.Range("G4").Value = "合夥‘"
$ws.Range("G5").Value = "合夥"
$ws.Range("H2").Value = "investment"
$ws.Range("H3").Value = "investment"
$ws.Range("H4").Value = "investment"
$ws.Range("H5").Value = "investment"
$ws.Range("I2").Value = "normal"
$ws.Range("I3").Value = "normal"
$ws.Range("I4").Value = "normal"
$ws.Range("I5").Value = "normal"
$ws.Range("J2").Value = "2011-11-21"
$ws.Range("J3").Value = "2011-11-21"
$ws.Range("J4").Value = "2011-11-21"
$ws.Range("J5").Value = "2011-11-21"
$ws.Range("K2").Value = "潘維剛"
$ws.Range("K3").Value = "潘維剛"
$ws.Range("K4").Value = "潘維剛"
$ws.Range("K5").Value = "潘維剛"
$ws.Range("L2").Value = 678
$ws.Range("L3").Value = 678
$ws.Range("L4").Value = 678
$ws.Range("L5").Value = 678
$ws.Range("M2").Value = "tmpcafb1"
$ws.Range("M3").Value = "tmpcafb1"
$ws.Range("M4").Value = "tmpcafb1"
$ws.Range("M5").Value = "tmpcafb1"
$ws.Range("N2").Value = 109
$ws.Range("N3").Value = 110
$ws.Range("N4").Value = 111
$ws.Range("N5").Value = 112
